$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.202.32"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Range('E2').Value = '  +4.68%  '
$ws.Range('D3').Value = "'2.272.86"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Range('E3').Value = '  +4.05%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'254.40"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = "'0.641"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').Value = "'72.29"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Range('E7').Value = '  +6.19%  '
$ws.Range('D8').Value = "'0.674"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Range('E8').Value = '  +17.88%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = "'40.39"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Range('E10').Value = '  +9.44%  '
$ws.Range('D11').Value = "'0.0981"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Range('E11').Value = '  +5.24%  '
$ws.Range('D12').Value = "'59.54"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').Value = "'7.59"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Range('E13').Value = '  +8.43%  '
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = "'2.620.26"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Range('E15').Value = '  +4.47%  '
$ws.Range('B16').Value = "'Chainlink"
$ws.Cells.Item(16,2).Style = "Normal"
$ws.Range('C16').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(16,3).Style = "Normal"
$ws.Range('D16').Value = "'14.97"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Range('E16').Value = '  +4.31%  '
$ws.Range('B17').Value = "'Polygon"
$ws.Cells.Item(17,2).Style = "Normal"
$ws.Range('C17').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(17,3).Style = "Normal"
$ws.Range('D17').Value = "'0.894"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('D18').Value = "'2.286.43"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Range('E18').Value = '  +3.86%  '
$ws.Range('D19').Value = "'43.181.72"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Range('E19').Value = '  +4.82%  '
$ws.Range('D20').Value = "'0.0₃0986"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('D21').Value = "'6.30"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('D22').Value = "'73.73"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('D23').Value = "'238.03"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Range('E23').Value = '  +2.43%  '
$ws.Range('D24').Value = "'2.12"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Range('E24').Value = '  +4.90%  '
$ws.Range('D25').Value = "'3.92"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').Value = "'11.83"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = "'2.48"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  +8.76%  '
$ws.Range('D31').Value = "'168.41"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = "'21.33"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('D33').Value = "'0.129"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Range('E33').Value = '  +9.87%  '
$ws.Range('D34').Value = "'6.12"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Range('E34').Value = '  +12.54%  '
$ws.Range('D35').Value = "'0.0782"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Range('E35').Value = '  +4.85%  '
$ws.Range('E36').Value = '  +2.68%  '
$ws.Range('D37').Value = "'29.19"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Range('E37').Value = '  +11.36%  '
$ws.Range('D38').Value = "'4.75"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('D39').Value = "'4.18"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').Value = "'0.0324"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Range('D41').Value = "'2.30"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Range('E41').Value = '  +5.20%  '
$ws.Range('D42').Value = "'5.95"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Range('E42').Value = '  +5.27%  '
$ws.Range('D43').Value = "'12.63"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Range('E43').Value = '  +3.59%  '
$ws.Range('D44').Value = "'64.62"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('D45').Value = "'5.02"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').Value = "'9.03"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Range('E47').Value = '  +5.36%  '
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('D49').Value = "'1.20"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('D50').Value = "'1.00"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('E51').Value = '  +5.25%  '
